# Auto-generated Excel COM-interop script
# Applies numeric updates to H:N columns across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
# per the target diff (Leve profit recalculation values).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 123
$ws.Range("I4").Value = 123
$ws.Range("K4").Value = 123
$ws.Range("M4").Value = -9
$ws.Range("H33").Value = 163.5
$ws.Range("I33").Value = 169.85
$ws.Range("K33").Value = 169.85
$ws.Range("M33").Value = 59.15000000000001
$ws.Range("H74").Value = 3500
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2064
$ws.Range("H77").Value = 3500
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10320
$ws.Range("H100").Value = 2111
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 2222
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 2222
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -3304
$ws.Range("H127").Value = 450
$ws.Range("I127").Value = 450
$ws.Range("K127").Value = 1350
$ws.Range("M127").Value = 3610

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1881.5
$ws.Range("I30").Value = 1881.5
$ws.Range("K30").Value = 1881.5
$ws.Range("M30").Value = -1731.5
$ws.Range("H122").Value = 1037
$ws.Range("I122").Value = 1255.5
$ws.Range("J122").Value = 600
$ws.Range("K122").Value = 3766.5
$ws.Range("L122").Value = 1800
$ws.Range("M122").Value = -1316.5
$ws.Range("N122").Value = -6700
$ws.Range("H132").Value = 6631.25
$ws.Range("I132").Value = 2262.5
$ws.Range("K132").Value = 6787.5
$ws.Range("M132").Value = -4257.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 39499
$ws.Range("J25").Value = 39499
$ws.Range("L25").Value = 39499
$ws.Range("N25").Value = -39969
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H86").Value = 3369.5
$ws.Range("I86").Value = 3000
$ws.Range("J86").Value = 3492.6667
$ws.Range("K86").Value = 3000
$ws.Range("L86").Value = 3492.6667
$ws.Range("M86").Value = -1877
$ws.Range("N86").Value = -5738.6667
$ws.Range("H89").Value = 3369.5
$ws.Range("I89").Value = 3000
$ws.Range("J89").Value = 3492.6667
$ws.Range("K89").Value = 15000
$ws.Range("L89").Value = 17463.3335
$ws.Range("M89").Value = -9384
$ws.Range("N89").Value = -28695.3335
$ws.Range("H105").Value = 1060.6875
$ws.Range("I105").Value = 910
$ws.Range("K105").Value = 910
$ws.Range("M105").Value = 837

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H7").Value = 26.166666
$ws.Range("I7").Value = 13.4
$ws.Range("K7").Value = 13.4
$ws.Range("M7").Value = 99.59999999999999
$ws.Range("H43").Value = 29172.428
$ws.Range("J43").Value = 29172.428
$ws.Range("L43").Value = 29172.428
$ws.Range("N43").Value = -29540.428
$ws.Range("H58").Value = 350
$ws.Range("I58").Value = 350
$ws.Range("K58").Value = 350
$ws.Range("M58").Value = -147
$ws.Range("H64").Value = 75000
$ws.Range("J64").Value = 75000
$ws.Range("L64").Value = 75000
$ws.Range("N64").Value = -75496
$ws.Range("H67").Value = 75000
$ws.Range("J67").Value = 75000
$ws.Range("L67").Value = 75000
$ws.Range("N67").Value = -76716
$ws.Range("H99").Value = 3427.8
$ws.Range("I99").Value = 2163
$ws.Range("K99").Value = 2163
$ws.Range("M99").Value = -665
$ws.Range("H101").Value = 29172.428
$ws.Range("J101").Value = 29172.428
$ws.Range("L101").Value = 29172.428
$ws.Range("N101").Value = -35662.428
$ws.Range("H126").Value = 3427.8
$ws.Range("I126").Value = 2163
$ws.Range("K126").Value = 6489
$ws.Range("M126").Value = -4019
$ws.Range("H132").Value = 6545
$ws.Range("I132").Value = 5151.8335
$ws.Range("K132").Value = 15455.5005
$ws.Range("M132").Value = -12925.5005
$ws.Range("H134").Value = 3816.3333
$ws.Range("I134").Value = 3816.3333
$ws.Range("K134").Value = 11448.9999
$ws.Range("M134").Value = -8913.999899999999
$ws.Range("H136").Value = 350
$ws.Range("I136").Value = 350
$ws.Range("K136").Value = 1050
$ws.Range("M136").Value = 1500
$ws.Range("H141").Value = 609582.6
$ws.Range("J141").Value = 609582.6
$ws.Range("L141").Value = 609582.6
$ws.Range("N141").Value = -619942.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H51").Value = 1000
$ws.Range("J51").Value = 1000
$ws.Range("L51").Value = 3000
$ws.Range("N51").Value = -3920
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 15000
$ws.Range("M70").Value = -14685
$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 15000
$ws.Range("M73").Value = -13908
$ws.Range("H132").Value = 1131.6666
$ws.Range("I132").Value = 947.5
$ws.Range("K132").Value = 8527.5
$ws.Range("M132").Value = -5997.5
$ws.Range("H133").Value = 3216.3333
$ws.Range("I133").Value = 3216.3333
$ws.Range("K133").Value = 9648.999899999999
$ws.Range("M133").Value = -4588.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13666667
$ws.Range("I11").Value = 13666667
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 13666667
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -13666528
$ws.Range("N11").ClearContents()
$ws.Range("H14").Value = 250000
$ws.Range("I14").Value = 250000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 250000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -249832
$ws.Range("N14").ClearContents()
$ws.Range("H24").Value = 15000
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15346
$ws.Range("H107").Value = 832.3333
$ws.Range("I107").Value = 832.3333
$ws.Range("K107").Value = 832.3333
$ws.Range("M107").Value = 1087.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 81.25
$ws.Range("I2").Value = 81.25
$ws.Range("K2").Value = 81.25
$ws.Range("M2").Value = 30.75
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H24").Value = 1000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 1000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 1000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -1686
$ws.Range("H25").Value = 2509500
$ws.Range("I25").Value = 2509500
$ws.Range("K25").Value = 2509500
$ws.Range("M25").Value = -2509270
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H93").Value = 1477.75
$ws.Range("I93").Value = 1477.75
$ws.Range("K93").Value = 1477.75
$ws.Range("M93").Value = -229.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H5").Value = 980
$ws.Range("J5").Value = 980
$ws.Range("L5").Value = 980
$ws.Range("N5").Value = -1204
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H124").Value = 93333
$ws.Range("J124").Value = 93333
$ws.Range("L124").Value = 93333
$ws.Range("N124").Value = -103153
